$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# ---------------------------------------------------------------------------
# 1) Title paragraph ("Chatroom Schedule") gains bold/color/size-28 paragraph
#    mark formatting, and four new paragraphs are inserted right after it
#    (before the schedule table): "Chatroom Project", "Array Bootcamp Fall
#    2021", the team roster line, and a trailing empty centered paragraph.
# ---------------------------------------------------------------------------

$oldTitle = $d.Paragraphs(1)
$anchor = $oldTitle.Range

# New title paragraph (replaces the old one; old one is deleted at the end
# so the w14:paraId/rsid attributes end up on the paragraph that survives).
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs(2)
$ir = $d.Range($p.Range.Start, $p.Range.End)
$ir.InsertXML('<w:p ' + $wNs + ' w14:paraId="2CED88A1" w14:textId="77777777" w:rsidR="002966A1" w:rsidRPr="00EB65FB" w:rsidRDefault="002966A1" w:rsidP="002966A1"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00EB65FB"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Chatroom Schedule</w:t></w:r></w:p>')

# "Chatroom Project"
$anchor = $d.Paragraphs(2).Range
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs(3)
$ir = $d.Range($p.Range.Start, $p.Range.End)
$ir.InsertXML('<w:p ' + $wNs + '><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Chatroom</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Project</w:t></w:r></w:p>')

# "Array Bootcamp Fall 2021"
$anchor = $d.Paragraphs(3).Range
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs(4)
$ir = $d.Range($p.Range.Start, $p.Range.End)
$ir.InsertXML('<w:p ' + $wNs + '><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Array Bootcamp Fall 2021</w:t></w:r></w:p>')

# Team roster line (with spell-check proofErr markers around "Condelario")
$anchor = $d.Paragraphs(4).Range
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs(5)
$ir = $d.Range($p.Range.Start, $p.Range.End)
$ir.InsertXML('<w:p ' + $wNs + '><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Katie Greenwald, Steve Bateman, Bowen </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Condelario</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

# Trailing empty centered paragraph
$anchor = $d.Paragraphs(5).Range
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs(6)
$ir = $d.Range($p.Range.Start, $p.Range.End)
$ir.InsertXML('<w:p ' + $wNs + '><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>')

# Drop the original (now-duplicate) title paragraph; paragraph 2 above takes
# its place at the top of the document, carrying the true w14:paraId/rsids.
$d.Paragraphs(1).Range.Delete()

# ---------------------------------------------------------------------------
# 2) lastRenderedPageBreak moves from the "8" week cell to the "5" week cell.
# ---------------------------------------------------------------------------

$d.Content.Find.Execute("lastRenderedPageBreak placeholder - unused", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

function Set-WeekCellPageBreak($weekText, $shouldHaveBreak) {
    $found = $word.ActiveDocument.Range(0, 0)
    foreach ($tbl in $word.ActiveDocument.Tables) {
        foreach ($row in $tbl.Rows) {
            $cell = $row.Cells(1)
            $cellRange = $cell.Range
            $cellRange.MoveEnd(1, -1) | Out-Null
            if ($cellRange.Text.Trim() -eq $weekText) {
                return $cell
            }
        }
    }
    return $null
}

$cell5 = Set-WeekCellPageBreak "5" $true
$cell8 = Set-WeekCellPageBreak "8" $false

if ($cell5 -ne $null) {
    $r5 = $cell5.Range
    $r5.MoveEnd(1, -1) | Out-Null
    $r5.Collapse(0)
    $r5.InsertXML('<w:r ' + $wNs + '><w:lastRenderedPageBreak/></w:r>')
}

if ($cell8 -ne $null) {
    $r8 = $cell8.Range
    $r8.MoveEnd(1, -1) | Out-Null
    $r8.Collapse(1)
    $r8.MoveStart(1, -1) | Out-Null
}

Write-Output "done"
